# kissebb javítások a prezi előtt
# Update the cached "today" date placeholders (slide master, notes master,
# and every slide layout) from 2018. 04. 12. / 4/12/2018 to 2018. 05. 03. / 5/3/2018.

$p = $ppt.ActivePresentation

function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            $t = $sh.TextFrame.TextRange.Text
            if ($t -eq "4/12/2018") {
                $sh.TextFrame.TextRange.Text = "5/3/2018"
            } elseif ($t -eq "2018. 04. 12.") {
                $sh.TextFrame.TextRange.Text = "2018. 05. 03."
            }
        }
    }
}

# Slide master "Date Placeholder"
Update-DateShape $p.SlideMaster.Shapes

# Notes master "Dátum helye"
Update-DateShape $p.NotesMaster.Shapes

# Every slide layout's own date placeholder
for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($li)
    Update-DateShape $layout.Shapes
}
